# Updated Voltage Design calculations
# Adds a "Design Side" mirror calculation block (columns F:H) next to the
# existing Measurement Side block (columns A:D), with a bold header label
# and two yellow-highlighted design input cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Design Side block header (bold) ---
$ws.Range("F1").Value = "Design Side"
$ws.Range("F1").Font.Bold = $true

# --- Row 2: ADC Ref / value / unit ---
$ws.Range("F2").Value = "ADC Ref"
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = "V"

# --- Row 3: Max Input level / value (highlighted) / unit ---
$ws.Range("F3").Value = "Max Input level"
$ws.Range("G3").Value = 485
$ws.Range("G3").Interior.Color = 65535
$ws.Range("H3").Value = "V"

# --- Row 4: SF / formula ---
$ws.Range("F4").Value = "SF"
$ws.Range("G4").Formula = "=CEILING(G3/G2,1)"

# --- Row 5: CF / formula ---
$ws.Range("F5").Value = "CF"
$ws.Range("G5").Formula = "=1/G4"

# --- Row 6: R1 / value (highlighted) / unit ---
$ws.Range("F6").Value = "R1"
$ws.Range("G6").Value = 450
$ws.Range("G6").Interior.Color = 65535
$ws.Range("H6").Value = "Kohm"

# --- Row 7: R2 / formula / unit ---
$ws.Range("F7").Value = "R2"
$ws.Range("G7").Formula = "=G5/(1-G5) *G6"
$ws.Range("H7").Value = "Kohm"

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 23.86
$ws.Columns("F").ColumnWidth = 19.3
$ws.Columns("G").ColumnWidth = 10.3

# --- Selection + page setup ---
$ws.Range("J7").Select() | Out-Null
$ws.PageSetup.Orientation = 1
